{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// This script applies four text edits to the document:\n//  1. Merge the split runs in the \"It turns out...\" paragraph (no text\n//     change, just removes the run split after \"game \" / \"starts\").\n//  2. Change \"maybe his parent is inside?\" -> \"maybe his parent's inside?\"\n//     (and merges the split runs into one run).\n//  3. Change \"Well, even if his parent is inside the store\" ->\n//     \"Well, even if his parent's inside the store\" (and merges the split\n//     runs into one run).\n//  4. Change \"?Greta (embarrassed embarrassed): I\u2026 uh\u2026\" ->\n//     \"?Greta (neutral embarrassed): I\u2026 uh\u2026\".\n\nconst edits = [\n  {\n    find:\n      \"It turns out that I didn\\u2019t need to set an alarm after all, as I wake up a solid two hours before the game starts. It should take twenty minutes or so to walk there, which leaves me an hour and a half to get out of the house. Child\\u2019s play for an expert at rushing outside in the morning.\",\n    replace:\n      \"It turns out that I didn\\u2019t need to set an alarm after all, as I wake up a solid two hours before the game starts. It should take twenty minutes or so to walk there, which leaves me an hour and a half to get out of the house. Child\\u2019s play for an expert at rushing outside in the morning.\",\n  },\n  {\n    find:\n      \"As I pass by the convenience store, I notice a small figure staring at me from across the street, a young boy who couldn\\u2019t be more than two or three years old. There seems to be nobody else around, but maybe his parent is inside?\",\n    replace:\n      \"As I pass by the convenience store, I notice a small figure staring at me from across the street, a young boy who couldn\\u2019t be more than two or three years old. There seems to be nobody else around, but maybe his parent\\u2019s inside?\",\n  },\n  {\n    find:\n      \"Well, even if his parent is inside the store, it\\u2019d probably be a good idea to make sure he doesn\\u2019t wander off or something.\",\n    replace:\n      \"Well, even if his parent\\u2019s inside the store, it\\u2019d probably be a good idea to make sure he doesn\\u2019t wander off or something.\",\n  },\n  {\n    find: \"?Greta (embarrassed embarrassed): I\\u2026 uh\\u2026\",\n    replace: \"?Greta (neutral embarrassed): I\\u2026 uh\\u2026\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of edits) {\n  const found = body.search(find, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Applies four text edits to the document:\n#  1. Merge the split runs in the \"It turns out...\" paragraph (no wording\n#     change, just removes the run split after \"game \" / \"starts\").\n#  2. Change \"maybe his parent is inside?\" -> \"maybe his parent's inside?\"\n#     (and merges the split runs into one run).\n#  3. Change \"Well, even if his parent is inside the store\" ->\n#     \"Well, even if his parent's inside the store\" (and merges the split\n#     runs into one run).\n#  4. Change \"?Greta (embarrassed embarrassed): I... uh...\" ->\n#     \"?Greta (neutral embarrassed): I... uh...\".\n#\n# Find.Execute's ReplaceWith/Replace arguments are used (rather than just\n# assigning Range.Text) because when the found text is identical to the\n# replacement text (edit #1 - only the run split changes, not the wording)\n# a plain Range.Text assignment is a content no-op and the runs stay split;\n# Find.Execute(..., Replace:=wdReplaceAll) always rewrites the matched\n# range into a single run.\n\n$wdReplaceAll = 2\n\n$apos = [char]0x2019\n$ellipsis = [char]0x2026\n\n$d = $word.ActiveDocument\n\n$edits = @(\n    @{\n        Find    = \"It turns out that I didn\" + $apos + \"t need to set an alarm after all, as I wake up a solid two hours before the game starts. It should take twenty minutes or so to walk there, which leaves me an hour and a half to get out of the house. Child\" + $apos + \"s play for an expert at rushing outside in the morning.\"\n        Replace = \"It turns out that I didn\" + $apos + \"t need to set an alarm after all, as I wake up a solid two hours before the game starts. It should take twenty minutes or so to walk there, which leaves me an hour and a half to get out of the house. Child\" + $apos + \"s play for an expert at rushing outside in the morning.\"\n    },\n    @{\n        Find    = \"As I pass by the convenience store, I notice a small figure staring at me from across the street, a young boy who couldn\" + $apos + \"t be more than two or three years old. There seems to be nobody else around, but maybe his parent is inside?\"\n        Replace = \"As I pass by the convenience store, I notice a small figure staring at me from across the street, a young boy who couldn\" + $apos + \"t be more than two or three years old. There seems to be nobody else around, but maybe his parent\" + $apos + \"s inside?\"\n    },\n    @{\n        Find    = \"Well, even if his parent is inside the store, it\" + $apos + \"d probably be a good idea to make sure he doesn\" + $apos + \"t wander off or something.\"\n        Replace = \"Well, even if his parent\" + $apos + \"s inside the store, it\" + $apos + \"d probably be a good idea to make sure he doesn\" + $apos + \"t wander off or something.\"\n    },\n    @{\n        Find    = \"?Greta (embarrassed embarrassed): I\" + $ellipsis + \" uh\" + $ellipsis\n        Replace = \"?Greta (neutral embarrassed): I\" + $ellipsis + \" uh\" + $ellipsis\n    }\n)\n\nforeach ($edit in $edits) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($edit.Find, $true, $false, $false, $false, $false, $true, 1, $false, $edit.Replace, $wdReplaceAll)\n    if (-not $found) {\n        Write-Output \"WARNING: not found -> $($edit.Find)\"\n    }\n}\n"}
